{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" answer cells in the practice\n// table with a new set of problems/answers (exact 1:1 text swap per cell,\n// matching the commit's regenerated worksheet data).\nconst replacements = [\n  [\"616\u00f74=154, 0\", \"827\u00f79=91, 8\"],\n  [\"656\u00f73=218, 2\", \"601\u00f72=300, 1\"],\n  [\"433\u00f78=54, 1\", \"260\u00f73=86, 2\"],\n  [\"398\u00f77=56, 6\", \"425\u00f73=141, 2\"],\n  [\"572\u00f73=190, 2\", \"612\u00f72=306, 0\"],\n  [\"486\u00f72=243, 0\", \"344\u00f74=86, 0\"],\n  [\"217\u00f75=43, 2\", \"819\u00f76=136, 3\"],\n  [\"428\u00f76=71, 2\", \"980\u00f72=490, 0\"],\n  [\"584\u00f72=292, 0\", \"315\u00f72=157, 1\"],\n  [\"131\u00f79=14, 5\", \"458\u00f75=91, 3\"],\n  [\"369\u00f78=46, 1\", \"497\u00f78=62, 1\"],\n  [\"352\u00f73=117, 1\", \"743\u00f75=148, 3\"],\n  [\"450\u00f76=75, 0\", \"789\u00f79=87, 6\"],\n  [\"968\u00f79=107, 5\", \"270\u00f74=67, 2\"],\n  [\"916\u00f79=101, 7\", \"176\u00f78=22, 0\"],\n  [\"442\u00f77=63, 1\", \"953\u00f72=476, 1\"],\n  [\"165\u00f73=55, 0\", \"364\u00f76=60, 4\"],\n  [\"531\u00f72=265, 1\", \"796\u00f77=113, 5\"],\n  [\"745\u00f72=372, 1\", \"347\u00f77=49, 4\"],\n  [\"521\u00f78=65, 1\", \"489\u00f74=122, 1\"],\n  [\"907\u00f74=226, 3\", \"815\u00f76=135, 5\"],\n  [\"718\u00f72=359, 0\", \"648\u00f73=216, 0\"],\n  [\"638\u00f75=127, 3\", \"914\u00f72=457, 0\"],\n  [\"518\u00f73=172, 2\", \"186\u00f76=31, 0\"],\n  [\"935\u00f79=103, 8\", \"586\u00f78=73, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at 9a8706d\n# Replace the 25 \"three-digit / one-digit\" division answer cells in the\n# practice table with the regenerated set of problems/answers. Each old\n# string is unique in the document, so a scoped Find/Replace per pair is\n# sufficient and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"616\u00f74=154, 0\"; Replace = \"827\u00f79=91, 8\" },\n    @{ Find = \"656\u00f73=218, 2\"; Replace = \"601\u00f72=300, 1\" },\n    @{ Find = \"433\u00f78=54, 1\"; Replace = \"260\u00f73=86, 2\" },\n    @{ Find = \"398\u00f77=56, 6\"; Replace = \"425\u00f73=141, 2\" },\n    @{ Find = \"572\u00f73=190, 2\"; Replace = \"612\u00f72=306, 0\" },\n    @{ Find = \"486\u00f72=243, 0\"; Replace = \"344\u00f74=86, 0\" },\n    @{ Find = \"217\u00f75=43, 2\"; Replace = \"819\u00f76=136, 3\" },\n    @{ Find = \"428\u00f76=71, 2\"; Replace = \"980\u00f72=490, 0\" },\n    @{ Find = \"584\u00f72=292, 0\"; Replace = \"315\u00f72=157, 1\" },\n    @{ Find = \"131\u00f79=14, 5\"; Replace = \"458\u00f75=91, 3\" },\n    @{ Find = \"369\u00f78=46, 1\"; Replace = \"497\u00f78=62, 1\" },\n    @{ Find = \"352\u00f73=117, 1\"; Replace = \"743\u00f75=148, 3\" },\n    @{ Find = \"450\u00f76=75, 0\"; Replace = \"789\u00f79=87, 6\" },\n    @{ Find = \"968\u00f79=107, 5\"; Replace = \"270\u00f74=67, 2\" },\n    @{ Find = \"916\u00f79=101, 7\"; Replace = \"176\u00f78=22, 0\" },\n    @{ Find = \"442\u00f77=63, 1\"; Replace = \"953\u00f72=476, 1\" },\n    @{ Find = \"165\u00f73=55, 0\"; Replace = \"364\u00f76=60, 4\" },\n    @{ Find = \"531\u00f72=265, 1\"; Replace = \"796\u00f77=113, 5\" },\n    @{ Find = \"745\u00f72=372, 1\"; Replace = \"347\u00f77=49, 4\" },\n    @{ Find = \"521\u00f78=65, 1\"; Replace = \"489\u00f74=122, 1\" },\n    @{ Find = \"907\u00f74=226, 3\"; Replace = \"815\u00f76=135, 5\" },\n    @{ Find = \"718\u00f72=359, 0\"; Replace = \"648\u00f73=216, 0\" },\n    @{ Find = \"638\u00f75=127, 3\"; Replace = \"914\u00f72=457, 0\" },\n    @{ Find = \"518\u00f73=172, 2\"; Replace = \"186\u00f76=31, 0\" },\n    @{ Find = \"935\u00f79=103, 8\"; Replace = \"586\u00f78=73, 2\" }\n)\n\nforeach ($item in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $null = $find.Execute(\n        $item.Find,    # FindText\n        $true,         # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $item.Replace, # ReplaceWith\n        2              # Replace (wdReplaceAll)\n    )\n}\n"}
